$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Update the Instance (D2) and TestDataFile (F2) values to reflect the new
# folder structure / added test scenarios.
$ws.Range("D2").Value = "Automation2"
$ws.Range("F2").Value = "WeeklySanity.xlsx"

# Update the active selection to D2 (as recorded in the saved view state).
$ws.Range("D2").Select()
